# Added Jira ID for RandomMerge Testcase
# Insert a new test-case row (ENWIAM55) right before the ENWIAM0003 block
# (old row 34), shifting the remaining rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Insert a blank row at 34 (old rows 34-37 become 35-38)
$ws.Rows.Item(34).Insert()

# Copy formatting (border + wrap-text) from the row just above (row 33,
# the ENWIAM54 row) so the new row matches the look of its neighbours.
$ws.Range("A33:D33").Copy($ws.Range("A34:D34"))

# New test-case content
$ws.Range("A34").Value = "ENWIAM55"
$ws.Range("B34").Value = "OPQA-2036"
$ws.Range("C34").Value = " From Neon, Verify that system is able to merge Activated STeAM account and Activated Facebook account and after merge verify STeAM TRUID is changed"
$ws.Range("D34").Value = "Y"

# Row height for the new row
$ws.Rows.Item(34).RowHeight = 45

# Restore the view/selection to mirror where the edit happened
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 29
$ws.Range("C34").Select() | Out-Null
